$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update Fitness (column C) values for the changed rows.
# Rows 2-12 (Generation 0-10): 7293 -> 7343
$ws.Range("C2:C12").Value = 7343

# Rows 13-18 (Generation 11-16): 7293 -> 7310
$ws.Range("C13:C18").Value = 7310
